$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 521.6087
$ws.Range("I19").Value = 589
$ws.Range("J19").Value = 469.76923
$ws.Range("K19").Value = 589
$ws.Range("L19").Value = 469.76923
$ws.Range("M19").Value = -414
$ws.Range("N19").Value = -819.76923
$ws.Range("H58").Value = 1184.9412
$ws.Range("I58").Value = 75.666664
$ws.Range("J58").Value = 1790
$ws.Range("K58").Value = 226.999992
$ws.Range("L58").Value = 5370
$ws.Range("M58").Value = -76.99999199999999
$ws.Range("N58").Value = -5670
$ws.Range("H69").Value = 3646.4348
$ws.Range("I69").Value = 2984.077
$ws.Range("J69").Value = 4507.5
$ws.Range("K69").Value = 8952.231
$ws.Range("L69").Value = 13522.5
$ws.Range("M69").Value = -8078.231
$ws.Range("N69").Value = -15270.5
$ws.Range("H72").Value = 3646.4348
$ws.Range("I72").Value = 2984.077
$ws.Range("J72").Value = 4507.5
$ws.Range("K72").Value = 26856.693
$ws.Range("L72").Value = 40567.5
$ws.Range("M72").Value = -22488.693
$ws.Range("N72").Value = -49303.5
$ws.Range("H74").Value = 3235.75
$ws.Range("I74").Value = 2817
$ws.Range("J74").Value = 3598.6667
$ws.Range("K74").Value = 2817
$ws.Range("L74").Value = 3598.6667
$ws.Range("M74").Value = -1881
$ws.Range("N74").Value = -5470.6667
$ws.Range("H77").Value = 3235.75
$ws.Range("I77").Value = 2817
$ws.Range("J77").Value = 3598.6667
$ws.Range("K77").Value = 14085
$ws.Range("L77").Value = 17993.3335
$ws.Range("M77").Value = -9405
$ws.Range("N77").Value = -27353.3335
$ws.Range("H100").Value = 22291610
$ws.Range("I100").Value = 37038930
$ws.Range("J100").Value = 170634.33
$ws.Range("K100").Value = 37038930
$ws.Range("L100").Value = 170634.33
$ws.Range("M100").Value = -37038389
$ws.Range("N100").Value = -171716.33
$ws.Range("H111").Value = 62503924
$ws.Range("I111").Value = 1813.5714
$ws.Range("J111").Value = 111116680
$ws.Range("K111").Value = 5440.7142
$ws.Range("L111").Value = 333350040
$ws.Range("M111").Value = -2373.7142
$ws.Range("N111").Value = -333356174
$ws.Range("H113").Value = 7524.926
$ws.Range("I113").Value = 2698.3125
$ws.Range("J113").Value = 14545.454
$ws.Range("K113").Value = 2698.3125
$ws.Range("L113").Value = 14545.454
$ws.Range("M113").Value = 555.6875
$ws.Range("N113").Value = -21053.454
$ws.Range("H129").Value = 799.2646999999999
$ws.Range("I129").Value = 478.75
$ws.Range("J129").Value = 897.88464
$ws.Range("K129").Value = 1436.25
$ws.Range("L129").Value = 2693.65392
$ws.Range("M129").Value = 3563.75
$ws.Range("N129").Value = -12693.65392
$ws.Range("H141").Value = 4061.6667
$ws.Range("I141").Value = 4061.6667
$ws.Range("K141").Value = 12185.0001
$ws.Range("M141").Value = -7005.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5920.1445
$ws.Range("I32").Value = 4411.0757
$ws.Range("K32").Value = 4411.0757
$ws.Range("M32").Value = -4124.0757
$ws.Range("H45").Value = 1351.0769
$ws.Range("I45").Value = 1340.909
$ws.Range("J45").Value = 1407
$ws.Range("K45").Value = 1340.909
$ws.Range("L45").Value = 1407
$ws.Range("M45").Value = -963.9090000000001
$ws.Range("N45").Value = -2161
$ws.Range("H74").Value = 5393
$ws.Range("I74").Value = 7477.9443
$ws.Range("J74").Value = 1981.2727
$ws.Range("K74").Value = 7477.9443
$ws.Range("L74").Value = 1981.2727
$ws.Range("M74").Value = -6603.9443
$ws.Range("N74").Value = -3729.2727
$ws.Range("H77").Value = 5393
$ws.Range("I77").Value = 7477.9443
$ws.Range("J77").Value = 1981.2727
$ws.Range("K77").Value = 37389.7215
$ws.Range("L77").Value = 9906.363499999999
$ws.Range("M77").Value = -33021.7215
$ws.Range("N77").Value = -18642.3635
$ws.Range("H97").Value = 719.2308
$ws.Range("I97").Value = 679.1667
$ws.Range("K97").Value = 679.1667
$ws.Range("M97").Value = -183.1667
$ws.Range("H110").Value = 128502.75
$ws.Range("I110").Value = 169670.33
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 169670.33
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -167625.33
$ws.Range("N110").Value = -9090

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 18535.166
$ws.Range("I20").Value = 14550.75
$ws.Range("J20").Value = 26504
$ws.Range("K20").Value = 14550.75
$ws.Range("L20").Value = 26504
$ws.Range("M20").Value = -14303.75
$ws.Range("N20").Value = -26998
$ws.Range("H80").Value = 2714.0667
$ws.Range("I80").Value = 875.7
$ws.Range("J80").Value = 3633.25
$ws.Range("K80").Value = 875.7
$ws.Range("L80").Value = 3633.25
$ws.Range("M80").Value = 122.3
$ws.Range("N80").Value = -5629.25
$ws.Range("H83").Value = 2714.0667
$ws.Range("I83").Value = 875.7
$ws.Range("J83").Value = 3633.25
$ws.Range("K83").Value = 4378.5
$ws.Range("L83").Value = 18166.25
$ws.Range("M83").Value = 613.5
$ws.Range("N83").Value = -28150.25
$ws.Range("H86").Value = 3479.818
$ws.Range("I86").Value = 3335
$ws.Range("J86").Value = 3866
$ws.Range("K86").Value = 3335
$ws.Range("L86").Value = 3866
$ws.Range("M86").Value = -2212
$ws.Range("N86").Value = -6112
$ws.Range("H89").Value = 3479.818
$ws.Range("I89").Value = 3335
$ws.Range("J89").Value = 3866
$ws.Range("K89").Value = 16675
$ws.Range("L89").Value = 19330
$ws.Range("M89").Value = -11059
$ws.Range("N89").Value = -30562
$ws.Range("H94").Value = 479.2857
$ws.Range("I94").Value = 399
$ws.Range("J94").Value = 680
$ws.Range("K94").Value = 399
$ws.Range("L94").Value = 680
$ws.Range("M94").Value = 52
$ws.Range("N94").Value = -1582
$ws.Range("H138").Value = 46156
$ws.Range("J138").Value = 46156
$ws.Range("L138").Value = 46156
$ws.Range("N138").Value = -56436

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3770.2778
$ws.Range("I62").Value = 2951.6667
$ws.Range("J62").Value = 4588.8887
$ws.Range("K62").Value = 2951.6667
$ws.Range("L62").Value = 4588.8887
$ws.Range("M62").Value = -2327.6667
$ws.Range("N62").Value = -5836.8887
$ws.Range("H65").Value = 3770.2778
$ws.Range("I65").Value = 2951.6667
$ws.Range("J65").Value = 4588.8887
$ws.Range("K65").Value = 14758.3335
$ws.Range("L65").Value = 22944.4435
$ws.Range("M65").Value = -11638.3335
$ws.Range("N65").Value = -29184.4435
$ws.Range("H109").Value = 27390
$ws.Range("J109").Value = 27390
$ws.Range("L109").Value = 27390
$ws.Range("N109").Value = -29470

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 820.0700000000001
$ws.Range("I131").Value = 575.8
$ws.Range("J131").Value = 832.92633
$ws.Range("K131").Value = 1727.4
$ws.Range("L131").Value = 2498.77899
$ws.Range("M131").Value = 3312.6
$ws.Range("N131").Value = -12578.77899
$ws.Range("H133").Value = 1575.238
$ws.Range("I133").Value = 1008.8889
$ws.Range("J133").Value = 2000
$ws.Range("K133").Value = 3026.6667
$ws.Range("L133").Value = 6000
$ws.Range("M133").Value = 2033.3333
$ws.Range("N133").Value = -16120

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2083.6667
$ws.Range("J41").Value = 5000
$ws.Range("L41").Value = 5000
$ws.Range("N41").Value = -5710
$ws.Range("H80").Value = 2087.5
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 2066.6667
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 2066.6667
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -4062.6667
$ws.Range("H83").Value = 2087.5
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 2066.6667
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 10333.3335
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -20317.3335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 84841.75
$ws.Range("I46").Value = 92463.73
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 92463.73
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -92275.73
$ws.Range("N46").Value = -1376
$ws.Range("H82").Value = 1510.6428
$ws.Range("I82").Value = 1437
$ws.Range("J82").Value = 1608.8334
$ws.Range("K82").Value = 1437
$ws.Range("L82").Value = 1608.8334
$ws.Range("M82").Value = -1076
$ws.Range("N82").Value = -2330.8334
$ws.Range("H85").Value = 1510.6428
$ws.Range("I85").Value = 1437
$ws.Range("J85").Value = 1608.8334
$ws.Range("K85").Value = 1437
$ws.Range("L85").Value = 1608.8334
$ws.Range("M85").Value = -189
$ws.Range("N85").Value = -4104.8334
$ws.Range("H93").Value = 8508
$ws.Range("I93").Value = 15020.857
$ws.Range("J93").Value = 909.6667
$ws.Range("K93").Value = 15020.857
$ws.Range("L93").Value = 909.6667
$ws.Range("M93").Value = -13772.857
$ws.Range("N93").Value = -3405.6667
$ws.Range("H132").Value = 5565.965
$ws.Range("I132").Value = 8054.278
$ws.Range("J132").Value = 3737.8164
$ws.Range("K132").Value = 24162.834
$ws.Range("L132").Value = 11213.4492
$ws.Range("M132").Value = -21632.834
$ws.Range("N132").Value = -16273.4492

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1357
$ws.Range("I107").Value = 1631.5555
$ws.Range("J107").Value = 533.3333
$ws.Range("K107").Value = 4894.666499999999
$ws.Range("L107").Value = 1599.9999
$ws.Range("M107").Value = -2974.666499999999
$ws.Range("N107").Value = -5439.9999
$ws.Range("H132").Value = 1163.4036
$ws.Range("I132").Value = 715.2273
$ws.Range("K132").Value = 2145.6819
$ws.Range("M132").Value = 384.3181
